# Insert two new data rows right after row 139 (i.e. at row 140), pushing the
# existing rows 140:243 down to 142:245. The new rows capture two additional
# "Acelga" price observations at Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 140 (Excel shifts 140:243 down to 142:245,
# copying the formatting of the row above into the new rows).
$ws.Rows("140:141").Insert()

# New row 140
$ws.Cells.Item(140, 1).Value = 7
$ws.Cells.Item(140, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(140, 3).Value = "Ñuble"
$ws.Cells.Item(140, 4).Value = 44762
$ws.Cells.Item(140, 5).Value = 16
$ws.Cells.Item(140, 6).Value = 100112009
$ws.Cells.Item(140, 7).Value = "Acelga"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 200
$ws.Cells.Item(140, 11).Value = 600
$ws.Cells.Item(140, 12).Value = 700
$ws.Cells.Item(140, 13).Value = 650
$ws.Cells.Item(140, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(140, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(140, 16).Value = 650
$ws.Cells.Item(140, 17).Value = 1
$ws.Cells.Item(140, 18).Value = "Hortaliza"

# New row 141
$ws.Cells.Item(141, 1).Value = 7
$ws.Cells.Item(141, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(141, 3).Value = "Ñuble"
$ws.Cells.Item(141, 4).Value = 44762
$ws.Cells.Item(141, 5).Value = 16
$ws.Cells.Item(141, 6).Value = 100112009
$ws.Cells.Item(141, 7).Value = "Acelga"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Segunda"
$ws.Cells.Item(141, 10).Value = 150
$ws.Cells.Item(141, 11).Value = 500
$ws.Cells.Item(141, 12).Value = 500
$ws.Cells.Item(141, 13).Value = 500
$ws.Cells.Item(141, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(141, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(141, 16).Value = 500
$ws.Cells.Item(141, 17).Value = 1
$ws.Cells.Item(141, 18).Value = "Hortaliza"
